$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New small helper cell at J2 (two spaces, preserved whitespace)
$ws.Range("J2").Value = "  "

# Row 7 now holds what used to be row 8's content (RANGE.GET)
$ws.Range("B7").Value = "RANGE.GET"
$ws.Range("C7").Value = "Return a range from a handle. Use cell menu 'Adjust' to get the entire range."

# Row 8 becomes the new JSON.PARSE entry
$ws.Range("B8").Value = "JSON.PARSE"
$ws.Range("C8").Value = "Parse a string into a key-value range. Keys starting with an asterisk (*) are handles to subranges."

# Row 9 (new row) gets the old JSON.GET entry, with a new description, copying B8's
# label formatting down since it is a fresh row
$ws.Range("B8").Copy($ws.Range("B9"))
$ws.Range("B9").Value = "JSON.GET"
$ws.Range("C9").Value = "Get JSON data from a parsed object."

# Row 10 (new row) gets the new JSON.VALUE entry
$ws.Range("B8").Copy($ws.Range("B10"))
$ws.Range("B10").Value = "JSON.VALUE"
$ws.Range("C10").Value = "Get values from a JSON object using keys."

# Selection ends up on C14 after the edits
$ws.Range("C14").Select()
